$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Monday" (column C) cell for these teacher rows - they were
# erroneously marked with an on-call tally of 1 on day 1; revert that.
$rows = @(9, 11, 12, 14, 28, 33, 47, 48, 49, 52, 56, 65)

foreach ($r in $rows) {
    $ws.Range("C$r").Clear()
}

$excel.CalculateFullRebuild()
